$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week range) ---
$ws.Range("A8").Value = "Volume 30   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/4/2023  Through  9/10/2023"

# --- Row 14 ---
$ws.Range("N14").Value = -50

# --- Row 15 ---
$ws.Range("N15").Value = -72.727272727272

# --- Row 16 (D16/E16 switch from blank-placeholder text to real numbers) ---
$ws.Range("D16").Value = 3
$ws.Range("D16").NumberFormat = "#,##0"
$ws.Range("E16").Value = -100
$ws.Range("E16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -42.857142857142
$ws.Range("J16").Value = 36
$ws.Range("K16").Value = 33.333333333333
$ws.Range("L16").Value = 37.142857142857
$ws.Range("M16").Value = -31.428571428571
$ws.Range("N16").Value = -87.5

# --- Row 17 ---
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = -50
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = -10
$ws.Range("I17").Value = 99
$ws.Range("J17").Value = 77
$ws.Range("K17").Value = 28.571428571428
$ws.Range("L17").Value = 59.677419354838
$ws.Range("M17").Value = 39.436619718309
$ws.Range("N17").Value = -51.470588235294

# --- Row 18 ---
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -57.142857142857
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -61.111111111111
$ws.Range("I18").Value = 70
$ws.Range("J18").Value = 69
$ws.Range("K18").Value = 1.449275362318
$ws.Range("L18").Value = -11.392405063291
$ws.Range("M18").Value = -59.537572254335
$ws.Range("N18").Value = -90.885416666666

# --- Row 19 ---
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -54.545454545454
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = -23.809523809523
$ws.Range("I19").Value = 300
$ws.Range("J19").Value = 374
$ws.Range("K19").Value = -19.786096256684
$ws.Range("L19").Value = 20.967741935483
$ws.Range("M19").Value = 38.888888888888
$ws.Range("N19").Value = -3.536977491961

# --- Row 20 (C20 switches from blank-placeholder text to a real number) ---
$ws.Range("C20").Value = 2
$ws.Range("C20").NumberFormat = "#,##0"
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 80
$ws.Range("I20").Value = 84
$ws.Range("J20").Value = 77
$ws.Range("K20").Value = 9.090909090909
$ws.Range("L20").Value = 61.538461538461
$ws.Range("M20").Value = -22.222222222222
$ws.Range("N20").Value = -93.745346239761

# --- Row 21 ---
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -54.166666666666
$ws.Range("F21").Value = 61
$ws.Range("G21").Value = 82
$ws.Range("H21").Value = -25.609756097561
$ws.Range("I21").Value = 609
$ws.Range("J21").Value = 642
$ws.Range("K21").Value = -5.140186915887
$ws.Range("L21").Value = 25.308641975308
$ws.Range("M21").Value = -5.581395348837
$ws.Range("N21").Value = -79.940711462450

# --- Row 24 ---
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -6.896551724137
$ws.Range("F24").Value = 129
$ws.Range("G24").Value = 150
$ws.Range("H24").Value = -14
$ws.Range("I24").Value = 1063
$ws.Range("J24").Value = 1284
$ws.Range("K24").Value = -17.211838006230
$ws.Range("L24").Value = 57.248520710059
$ws.Range("M24").Value = 32.213930348258

# --- Row 25 ---
$ws.Range("C25").Value = 5
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 31
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = 82.352941176470
$ws.Range("I25").Value = 262
$ws.Range("J25").Value = 220
$ws.Range("K25").Value = 19.090909090909
$ws.Range("L25").Value = 52.325581395348
$ws.Range("M25").Value = 1.158301158301

# --- Row 27 ---
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 22
$ws.Range("K27").Value = -33.333333333333
$ws.Range("L27").Value = 29.411764705882

# --- Row 28 (C28 switches from a real number back to the blank-placeholder text) ---
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("D28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("N28").Value = -20

# --- Row 29 (C29 switches from a real number back to the blank-placeholder text) ---
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("D29").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("N29").Value = -25
